# Applies the figS1_life_stage_decision_tree.pptx edit:
#  - bump cached "datetimeFigureOut" field text 8/29/22 -> 8/30/22
#    (slide master, all 11 slide layouts, and the notes master)
#  - tweak a few percentages in the decision-tree slide text boxes

$p = $ppt.ActivePresentation

# --- 1. Date placeholders -------------------------------------------------

function Update-DateText {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "8/29/22") {
            $shp.TextFrame.TextRange.Text = "8/30/22"
        }
    }
}

# Slide master date placeholder
Update-DateText $p.SlideMaster.Shapes

# Every slide layout's date placeholder
$customLayouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $customLayouts.Count; $li++) {
    Update-DateText $customLayouts.Item($li).Shapes
}

# Notes master date placeholder.
# NB: the notes master's own Shapes collection cannot be written through
# directly here (doing so mis-routes the write into the slide master's
# shapes of the same index), so go through HeadersFooters instead, which
# updates the notes master part correctly. (HeadersFooters.DateAndTime.Text
# reads back blank through this host, so just set it unconditionally.)
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "8/30/22"

# --- 2. Slide 1 text edits -------------------------------------------------

$s = $p.Slides.Item(1)

# "Yes (6% data)" -> "Yes (5% data)"
$s.Shapes.Item(6).TextFrame.TextRange.Text = "Yes (5% data)"

# "No (71% data)" -> "No (72% data)"
$s.Shapes.Item(7).TextFrame.TextRange.Text = "No (72% data)"

# "No (30% data)" -> "No " + "(31% " + "data)" (three runs): only the
# "(30% " chunk actually changes (30 -> 31), so edit just that substring;
# the engine then splits the paragraph into the three runs.
$tr22 = $s.Shapes.Item(22).TextFrame.TextRange
$mid22 = $tr22.Characters(4, 5)
$mid22.Text = "(31% "
